$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the investor names in column A (rows 2-7) with generic placeholders,
# as part of adding tests for fund uploads.
$ws.Range("A2").Value = "Investor 1"
$ws.Range("A3").Value = "Investor 2"
$ws.Range("A4").Value = "Investor 3"
$ws.Range("A5").Value = "Investor 4"
$ws.Range("A6").Value = "Investor 5"
$ws.Range("A7").Value = "Investor 6"

# Apply a distinct cell style to the updated investor-name cells.
$st = $wb.Styles.Add("Normal 3")
$st.Font.Name = "Arial"
$st.Font.Size = 11
$ws.Range("A2:A7").Style = "Normal 3"

# Update the active selection to reflect the edited column.
$ws.Range("A2:A6").Select() | Out-Null
